$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B18 was stored as an inline string "4"; convert it to a real number 4
$ws.Range("B18").Value = 4

# Add new row 19 with annotation data
$ws.Range("A19").Value = "Ying Tang"

# B19 stays a text value "5" (not numeric). Enter it as a text formula then
# paste-special as values so it lands as a literal string without requiring
# any NumberFormat change (which would otherwise mint an unused style).
$ws.Range("B19").Formula = '="5"'
$ws.Range("B19").Copy()
$ws.Range("B19").PasteSpecial(-4163)

$ws.Range("C19").Value = "thank everybody again,useful suggestions"
$ws.Range("D19").Value = "ACK"
$ws.Range("E19").Value = "OTH"
$ws.Range("F19").Value = "e27c53be-a9c3-4697-b8f8-f90bcc73c090"
$ws.Range("G19").Value = "SJaP_-xAb_annotated.xlsx"
$ws.Range("H19").Value = "We thank everybody again for their useful suggestions and we uploaded a revision of the paper."
